$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3) appended below the existing header (row 1) and
# single data row (row 2).

# --- Plain numeric cells ---
$ws.Range("A3").Value = 112111626
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("Q3").Value = 555846.0651465225
$ws.Range("R3").Value = 6952042.273423757
$ws.Range("S3").Value = 10

# --- Plain text cells ---
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("M3").Value = "födosökande"
$ws.Range("N3").Value = "observerad"
$ws.Range("P3").Value = "Andersloken, Mpd"
$ws.Range("T3").Value = "Västernorrland"
$ws.Range("U3").Value = "Ånge"
$ws.Range("V3").Value = "Medelpad"
$ws.Range("W3").Value = "Borgsjö"
$ws.Range("AW3").Value = "Benny Öwre"
$ws.Range("AX3").Value = "Benny Öwre"

# --- Text cells whose content looks numeric / date-like: force text with
#     a leading apostrophe so they are not auto-coerced to number/date ---
$ws.Range("I3").Value = "'1"
$ws.Range("Y3").Value = "'2023-09-14"
$ws.Range("Z3").Value = "'00:00"
$ws.Range("AA3").Value = "'2023-09-14"
$ws.Range("AB3").Value = "'00:00"

# --- Explicitly-present but empty text cells ---
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("AT3").Value = "'"
$ws.Range("AY3").Value = "'"

# --- Boolean cells ---
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
